$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(141, 8).Value = 3064.3333
$ws.Cells.Item(141, 9).Value = 2252.3157
$ws.Cells.Item(141, 11).Value = 6756.9471
$ws.Cells.Item(141, 13).Value = -1576.9471

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1999.5
$ws.Cells.Item(45, 9).Value = 1999
$ws.Cells.Item(45, 11).Value = 1999
$ws.Cells.Item(45, 13).Value = -1622
$ws.Cells.Item(46, 8).Value = 26666.666
$ws.Cells.Item(46, 9).Value = 26666.666
$ws.Cells.Item(46, 11).Value = 26666.666
$ws.Cells.Item(46, 13).Value = -26347.666
$ws.Cells.Item(61, 8).Value = 2855.125
$ws.Cells.Item(61, 9).Value = 3144.3333
$ws.Cells.Item(61, 11).Value = 3144.3333
$ws.Cells.Item(61, 13).Value = -2932.3333
$ws.Cells.Item(80, 8).Value = 50073
$ws.Cells.Item(80, 10).Value = 50073
$ws.Cells.Item(80, 12).Value = 50073
$ws.Cells.Item(80, 14).Value = -52069
$ws.Cells.Item(83, 8).Value = 50073
$ws.Cells.Item(83, 10).Value = 50073
$ws.Cells.Item(83, 12).Value = 150219
$ws.Cells.Item(83, 14).Value = -160203
$ws.Cells.Item(124, 8).Value = 37885.6
$ws.Cells.Item(124, 10).Value = 37885.6
$ws.Cells.Item(124, 12).Value = 37885.6
$ws.Cells.Item(124, 14).Value = -47705.6
$ws.Cells.Item(129, 8).Value = 78000
$ws.Cells.Item(129, 10).Value = 78000
$ws.Cells.Item(129, 12).Value = 78000
$ws.Cells.Item(129, 14).Value = -88000
$ws.Cells.Item(132, 8).Value = 1196.6923
$ws.Cells.Item(132, 9).Value = 1214.3636
$ws.Cells.Item(132, 10).Value = 1099.5
$ws.Cells.Item(132, 11).Value = 3643.0908
$ws.Cells.Item(132, 12).Value = 3298.5
$ws.Cells.Item(132, 13).Value = -1113.0908
$ws.Cells.Item(132, 14).Value = -8358.5
$ws.Cells.Item(136, 8).Value = 2855.125
$ws.Cells.Item(136, 9).Value = 3144.3333
$ws.Cells.Item(136, 11).Value = 9432.999899999999
$ws.Cells.Item(136, 13).Value = -6882.999899999999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(75, 8).Value = 48749
$ws.Cells.Item(75, 10).Value = 59998.668
$ws.Cells.Item(75, 12).Value = 59998.668
$ws.Cells.Item(75, 14).Value = -61870.668
$ws.Cells.Item(78, 8).Value = 48749
$ws.Cells.Item(78, 10).Value = 59998.668
$ws.Cells.Item(78, 12).Value = 179996.004
$ws.Cells.Item(78, 14).Value = -189356.004
$ws.Cells.Item(80, 8).Value = 498.66666
$ws.Cells.Item(80, 9).Value = 574.4
$ws.Cells.Item(80, 10).Value = 404
$ws.Cells.Item(80, 11).Value = 574.4
$ws.Cells.Item(80, 12).Value = 404
$ws.Cells.Item(80, 13).Value = 423.6
$ws.Cells.Item(80, 14).Value = -2400
$ws.Cells.Item(82, 14).ClearContents()
$ws.Cells.Item(82, 8).Value = 27500
$ws.Cells.Item(82, 9).Value = 27500
$ws.Cells.Item(82, 10).Value = 0
$ws.Cells.Item(82, 11).Value = 27500
$ws.Cells.Item(82, 12).Value = 0
$ws.Cells.Item(82, 13).Value = -27117
$ws.Cells.Item(83, 8).Value = 498.66666
$ws.Cells.Item(83, 9).Value = 574.4
$ws.Cells.Item(83, 10).Value = 404
$ws.Cells.Item(83, 11).Value = 2872
$ws.Cells.Item(83, 12).Value = 2020
$ws.Cells.Item(83, 13).Value = 2120
$ws.Cells.Item(83, 14).Value = -12004
$ws.Cells.Item(85, 14).ClearContents()
$ws.Cells.Item(85, 8).Value = 27500
$ws.Cells.Item(85, 9).Value = 27500
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 11).Value = 27500
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 13).Value = -26174
$ws.Cells.Item(99, 8).Value = 34746.934
$ws.Cells.Item(99, 9).Value = 42962.832
$ws.Cells.Item(99, 10).Value = 1883.3334
$ws.Cells.Item(99, 11).Value = 42962.832
$ws.Cells.Item(99, 12).Value = 1883.3334
$ws.Cells.Item(99, 13).Value = -41464.832
$ws.Cells.Item(99, 14).Value = -4879.3334

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 14).ClearContents()
$ws.Cells.Item(22, 8).Value = 241
$ws.Cells.Item(22, 9).Value = 241
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 241
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = 109
$ws.Cells.Item(31, 8).Value = 3310.6487
$ws.Cells.Item(31, 9).Value = 3234.4
$ws.Cells.Item(31, 10).Value = 3362.6365
$ws.Cells.Item(31, 11).Value = 3234.4
$ws.Cells.Item(31, 12).Value = 3362.6365
$ws.Cells.Item(31, 13).Value = -2939.4
$ws.Cells.Item(31, 14).Value = -3952.6365
$ws.Cells.Item(34, 8).Value = 3310.6487
$ws.Cells.Item(34, 9).Value = 3234.4
$ws.Cells.Item(34, 10).Value = 3362.6365
$ws.Cells.Item(34, 11).Value = 3234.4
$ws.Cells.Item(34, 12).Value = 3362.6365
$ws.Cells.Item(34, 13).Value = -3032.4
$ws.Cells.Item(34, 14).Value = -3766.6365
$ws.Cells.Item(99, 8).Value = 6743.643
$ws.Cells.Item(99, 9).Value = 2648.875
$ws.Cells.Item(99, 11).Value = 2648.875
$ws.Cells.Item(99, 13).Value = -1150.875
$ws.Cells.Item(126, 8).Value = 6743.643
$ws.Cells.Item(126, 9).Value = 2648.875
$ws.Cells.Item(126, 11).Value = 7946.625
$ws.Cells.Item(126, 13).Value = -5476.625
$ws.Cells.Item(132, 8).Value = 2689.4167
$ws.Cells.Item(132, 9).Value = 2025
$ws.Cells.Item(132, 11).Value = 6075
$ws.Cells.Item(132, 13).Value = -3545
$ws.Cells.Item(134, 8).Value = 1781.7646
$ws.Cells.Item(134, 9).Value = 1536.7368
$ws.Cells.Item(134, 11).Value = 4610.2104
$ws.Cells.Item(134, 13).Value = -2075.2104

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 801
$ws.Cells.Item(68, 9).Value = 602
$ws.Cells.Item(68, 11).Value = 1806
$ws.Cells.Item(68, 13).Value = -995
$ws.Cells.Item(71, 8).Value = 801
$ws.Cells.Item(71, 9).Value = 602
$ws.Cells.Item(71, 11).Value = 5418
$ws.Cells.Item(71, 13).Value = -1362
$ws.Cells.Item(80, 14).ClearContents()
$ws.Cells.Item(80, 8).Value = 1942.25
$ws.Cells.Item(80, 9).Value = 1942.25
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 11).Value = 5826.75
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 13).Value = -4890.75
$ws.Cells.Item(83, 14).ClearContents()
$ws.Cells.Item(83, 8).Value = 1942.25
$ws.Cells.Item(83, 9).Value = 1942.25
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 11).Value = 17480.25
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 13).Value = -12800.25
$ws.Cells.Item(86, 8).Value = 375
$ws.Cells.Item(86, 9).Value = 250
$ws.Cells.Item(86, 11).Value = 750
$ws.Cells.Item(86, 13).Value = 436
$ws.Cells.Item(89, 8).Value = 375
$ws.Cells.Item(89, 9).Value = 250
$ws.Cells.Item(89, 11).Value = 2250
$ws.Cells.Item(89, 13).Value = 3678
$ws.Cells.Item(129, 8).Value = 2837.7
$ws.Cells.Item(129, 9).Value = 2194.5
$ws.Cells.Item(129, 10).Value = 2998.5
$ws.Cells.Item(129, 11).Value = 6583.5
$ws.Cells.Item(129, 12).Value = 8995.5
$ws.Cells.Item(129, 13).Value = -1583.5
$ws.Cells.Item(129, 14).Value = -18995.5
$ws.Cells.Item(130, 8).Value = 5411.75
$ws.Cells.Item(130, 9).Value = 1323.5
$ws.Cells.Item(130, 10).Value = 9500
$ws.Cells.Item(130, 11).Value = 3970.5
$ws.Cells.Item(130, 12).Value = 28500
$ws.Cells.Item(130, 13).Value = 1049.5
$ws.Cells.Item(130, 14).Value = -38540

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2446.2856
$ws.Cells.Item(80, 9).Value = 1924.8
$ws.Cells.Item(80, 10).Value = 3750
$ws.Cells.Item(80, 11).Value = 1924.8
$ws.Cells.Item(80, 12).Value = 3750
$ws.Cells.Item(80, 13).Value = -926.8
$ws.Cells.Item(80, 14).Value = -5746
$ws.Cells.Item(83, 8).Value = 2446.2856
$ws.Cells.Item(83, 9).Value = 1924.8
$ws.Cells.Item(83, 10).Value = 3750
$ws.Cells.Item(83, 11).Value = 9624
$ws.Cells.Item(83, 12).Value = 18750
$ws.Cells.Item(83, 13).Value = -4632
$ws.Cells.Item(83, 14).Value = -28734
$ws.Cells.Item(132, 8).Value = 2848.35
$ws.Cells.Item(132, 9).Value = 2763.2104
$ws.Cells.Item(132, 11).Value = 8289.6312
$ws.Cells.Item(132, 13).Value = -5759.6312

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3862.5557
$ws.Cells.Item(7, 9).Value = 3758.7144
$ws.Cells.Item(7, 10).Value = 4226
$ws.Cells.Item(7, 11).Value = 3758.7144
$ws.Cells.Item(7, 12).Value = 4226
$ws.Cells.Item(7, 13).Value = -3646.7144
$ws.Cells.Item(7, 14).Value = -4450
$ws.Cells.Item(16, 8).Value = 19424.625
$ws.Cells.Item(16, 9).Value = 19499.25
$ws.Cells.Item(16, 11).Value = 19499.25
$ws.Cells.Item(16, 13).Value = -19329.25
$ws.Cells.Item(82, 8).Value = 1818.125
$ws.Cells.Item(82, 10).Value = 2200
$ws.Cells.Item(82, 12).Value = 2200
$ws.Cells.Item(82, 14).Value = -2922
$ws.Cells.Item(85, 8).Value = 1818.125
$ws.Cells.Item(85, 10).Value = 2200
$ws.Cells.Item(85, 12).Value = 2200
$ws.Cells.Item(85, 14).Value = -4696
$ws.Cells.Item(126, 8).Value = 3862.5557
$ws.Cells.Item(126, 9).Value = 3758.7144
$ws.Cells.Item(126, 10).Value = 4226
$ws.Cells.Item(126, 11).Value = 11276.1432
$ws.Cells.Item(126, 12).Value = 12678
$ws.Cells.Item(126, 13).Value = -8806.143199999999
$ws.Cells.Item(126, 14).Value = -17618
$ws.Cells.Item(132, 8).Value = 4800.4614
$ws.Cells.Item(132, 9).Value = 3942.0625
$ws.Cells.Item(132, 10).Value = 6173.9
$ws.Cells.Item(132, 11).Value = 11826.1875
$ws.Cells.Item(132, 12).Value = 18521.7
$ws.Cells.Item(132, 13).Value = -9296.1875
$ws.Cells.Item(132, 14).Value = -23581.7
$ws.Cells.Item(136, 8).Value = 2758.0454
$ws.Cells.Item(136, 9).Value = 2552.1177
$ws.Cells.Item(136, 11).Value = 7656.353099999999
$ws.Cells.Item(136, 13).Value = -5106.353099999999

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 1003834
$ws.Cells.Item(5, 10).Value = 5751
$ws.Cells.Item(5, 12).Value = 5751
$ws.Cells.Item(5, 14).Value = -5975
$ws.Cells.Item(113, 14).ClearContents()
$ws.Cells.Item(113, 8).Value = 710.25
$ws.Cells.Item(113, 9).Value = 710.25
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 2130.75
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = 39.25
$ws.Cells.Item(126, 8).Value = 2290
$ws.Cells.Item(126, 9).Value = 2021.1111
$ws.Cells.Item(126, 11).Value = 6063.3333
$ws.Cells.Item(126, 13).Value = -3593.3333
$ws.Cells.Item(132, 8).Value = 32596.133
$ws.Cells.Item(132, 9).Value = 42019.22
$ws.Cells.Item(132, 11).Value = 126057.66
$ws.Cells.Item(132, 13).Value = -123527.66
$ws.Cells.Item(136, 8).Value = 2013.5
$ws.Cells.Item(136, 9).Value = 2017.4546
$ws.Cells.Item(136, 11).Value = 6052.3638
$ws.Cells.Item(136, 13).Value = -3502.3638
